$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows: row number, new Timestamp (col A) serial value, new Notified Production (col B) value
$data = @(
    @(2, 45736.01041666666, 0),
    @(3, 45736.02083333334, 0),
    @(4, 45736.03125, 0),
    @(5, 45736.04166666666, 0),
    @(6, 45736.05208333334, 0),
    @(7, 45736.0625, 0),
    @(8, 45736.07291666666, 0),
    @(9, 45736.08333333334, 0),
    @(10, 45736.09375, 0),
    @(11, 45736.10416666666, 0),
    @(12, 45736.11458333334, 0),
    @(13, 45736.125, 0),
    @(14, 45736.13541666666, 0),
    @(15, 45736.14583333334, 0),
    @(16, 45736.15625, 0),
    @(17, 45736.16666666666, 0),
    @(18, 45736.17708333334, 1),
    @(19, 45736.1875, 1),
    @(20, 45736.19791666666, 1),
    @(21, 45736.20833333334, 1),
    @(22, 45736.21875, 35),
    @(23, 45736.22916666666, 24),
    @(24, 45736.23958333334, 25),
    @(25, 45736.25, 35),
    @(26, 45736.26041666666, 278),
    @(27, 45736.27083333334, 309),
    @(28, 45736.28125, 345),
    @(29, 45736.29166666666, 402),
    @(30, 45736.30208333334, 895),
    @(31, 45736.3125, 948),
    @(32, 45736.32291666666, 1012),
    @(33, 45736.33333333334, 1074),
    @(34, 45736.34375, 1546),
    @(35, 45736.35416666666, 1597),
    @(36, 45736.36458333334, 1646),
    @(37, 45736.375, 1694),
    @(38, 45736.38541666666, 1989),
    @(39, 45736.39583333334, 2016),
    @(40, 45736.40625, 2039),
    @(41, 45736.41666666666, 2059),
    @(42, 45736.42708333334, 2175),
    @(43, 45736.4375, 2187),
    @(44, 45736.44791666666, 2197),
    @(45, 45736.45833333334, 2204),
    @(46, 45736.46875, 2226),
    @(47, 45736.47916666666, 2227),
    @(48, 45736.48958333334, 2227),
    @(49, 45736.5, 2224),
    @(50, 45736.51041666666, 2173),
    @(51, 45736.52083333334, 2165),
    @(52, 45736.53125, 2153),
    @(53, 45736.54166666666, 2136),
    @(54, 45736.55208333334, 2008),
    @(55, 45736.5625, 1984),
    @(56, 45736.57291666666, 1959),
    @(57, 45736.58333333334, 1931),
    @(58, 45736.59375, 1679),
    @(59, 45736.60416666666, 1640),
    @(60, 45736.61458333334, 1595),
    @(61, 45736.625, 1548),
    @(62, 45736.63541666666, 1061),
    @(63, 45736.64583333334, 1008),
    @(64, 45736.65625, 952),
    @(65, 45736.66666666666, 909),
    @(66, 45736.67708333334, 411),
    @(67, 45736.6875, 361),
    @(68, 45736.69791666666, 321),
    @(69, 45736.70833333334, 292),
    @(70, 45736.71875, 34),
    @(71, 45736.72916666666, 21),
    @(72, 45736.73958333334, 17),
    @(73, 45736.75, 16),
    @(74, 45736.76041666666, 2),
    @(75, 45736.77083333334, 2),
    @(76, 45736.78125, 2),
    @(77, 45736.79166666666, 2),
    @(78, 45736.80208333334, 2),
    @(79, 45736.8125, 2),
    @(80, 45736.82291666666, 2),
    @(81, 45736.83333333334, 2),
    @(82, 45736.84375, 1),
    @(83, 45736.85416666666, 0),
    @(84, 45736.86458333334, 0),
    @(85, 45736.875, 0),
    @(86, 45736.88541666666, 0),
    @(87, 45736.89583333334, 0),
    @(88, 45736.90625, 0),
    @(89, 45736.91666666666, 0),
    @(90, 45736.92708333334, 0),
    @(91, 45736.9375, 0),
    @(92, 45736.94791666666, 0),
    @(93, 45736.95833333334, 0),
    @(94, 45736.96875, 0),
    @(95, 45736.97916666666, 0),
    @(96, 45736.98958333334, 0),
    @(97, 45737, 0)
)

foreach ($item in $data) {
    $row = $item[0]
    $tsValue = $item[1]
    $prodValue = $item[2]
    $ws.Cells.Item($row, 1).Value = $tsValue
    $ws.Cells.Item($row, 2).Value = $prodValue
}
